$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add headers I1 and J1, copying the formatting (style) from H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for columns I (I0) and J (IF), rows 2-41
$data = @(
    @(9, 9),
    @(9, 9),
    @(8, 9),
    @(8, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(9, 9),
    @(8, 9),
    @(9, 9),
    @(9, 10),
    @(8, 9),
    @(8, 9),
    @(8, 9),
    @(8, 9),
    @(9, 9),
    @(9, 9),
    @(7, 8),
    @(8, 8),
    @(9, 9),
    @(7, 8),
    @(8, 9),
    @(9, 9),
    @(9, 9),
    @(8, 8),
    @(8, 9),
    @(9, 9),
    @(5, 6),
    @(6, 6),
    @(7, 7),
    @(7, 7),
    @(5, 5),
    @(4, 4),
    @(5, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
